$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 10245.77
$ws.Range("I28").Value = 5687.375
$ws.Range("J28").Value = 17539.2
$ws.Range("K28").Value = 5687.375
$ws.Range("L28").Value = 17539.2
$ws.Range("M28").Value = -5202.375
$ws.Range("N28").Value = -18509.2
$ws.Range("H58").Value = 14176.556
$ws.Range("J58").Value = 37500
$ws.Range("L58").Value = 112500
$ws.Range("N58").Value = -112800
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = 0
$ws.Range("H132").Value = 1966.9429
$ws.Range("I132").Value = 1373.4762
$ws.Range("K132").Value = 4120.4286
$ws.Range("M132").Value = -1590.4286

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6495810
$ws.Range("I2").Value = 6495810
$ws.Range("K2").Value = 6495810
$ws.Range("M2").Value = -6495697
$ws.Range("H102").Value = 2347.76
$ws.Range("I102").Value = 2347.76
$ws.Range("K102").Value = 2347.76
$ws.Range("M102").Value = -725.7600000000002
$ws.Range("H110").Value = 45001556
$ws.Range("I110").Value = 56251570
$ws.Range("K110").Value = 56251570
$ws.Range("M110").Value = -56249525
$ws.Range("H116").Value = 6495810
$ws.Range("I116").Value = 6495810
$ws.Range("K116").Value = 6495810
$ws.Range("M116").Value = -6493516
$ws.Range("H122").Value = 3012
$ws.Range("I122").Value = 3012
$ws.Range("K122").Value = 9036
$ws.Range("M122").Value = -6586

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6495810
$ws.Range("I3").Value = 6495810
$ws.Range("K3").Value = 6495810
$ws.Range("M3").Value = -6495696
$ws.Range("H20").Value = 4063.9375
$ws.Range("I20").Value = 3475.6
$ws.Range("J20").Value = 5044.5
$ws.Range("K20").Value = 3475.6
$ws.Range("L20").Value = 5044.5
$ws.Range("M20").Value = -3228.6
$ws.Range("N20").Value = -5538.5
$ws.Range("H29").Value = 94.333336
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H64").Value = 2543.889
$ws.Range("J64").Value = 3478.3333
$ws.Range("L64").Value = 3478.3333
$ws.Range("N64").Value = -3928.3333
$ws.Range("H67").Value = 2543.889
$ws.Range("J67").Value = 3478.3333
$ws.Range("L67").Value = 3478.3333
$ws.Range("N67").Value = -5038.3333

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4104.227
$ws.Range("I16").Value = 2627
$ws.Range("K16").Value = 2627
$ws.Range("M16").Value = -2340
$ws.Range("H31").Value = 62503644
$ws.Range("I31").Value = 71431490
$ws.Range("J31").Value = 8749.5
$ws.Range("K31").Value = 71431490
$ws.Range("L31").Value = 8749.5
$ws.Range("M31").Value = -71431195
$ws.Range("N31").Value = -9339.5
$ws.Range("H34").Value = 62503644
$ws.Range("I34").Value = 71431490
$ws.Range("J34").Value = 8749.5
$ws.Range("K34").Value = 71431490
$ws.Range("L34").Value = 8749.5
$ws.Range("M34").Value = -71431288
$ws.Range("N34").Value = -9153.5
$ws.Range("H99").Value = 2500
$ws.Range("I99").Value = 2500
$ws.Range("K99").Value = 2500
$ws.Range("M99").Value = -1002
$ws.Range("H113").Value = 4104.227
$ws.Range("I113").Value = 2627
$ws.Range("K113").Value = 2627
$ws.Range("M113").Value = -457
$ws.Range("H126").Value = 2500
$ws.Range("I126").Value = 2500
$ws.Range("K126").Value = 7500
$ws.Range("M126").Value = -5030
$ws.Range("H134").Value = 10741.143
$ws.Range("I134").Value = 10448.154
$ws.Range("K134").Value = 31344.462
$ws.Range("M134").Value = -28809.462

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 660.8
$ws.Range("I8").Value = 660.8
$ws.Range("K8").Value = 1982.4
$ws.Range("M8").Value = -1843.4
$ws.Range("H61").Value = 83.2
$ws.Range("J61").Value = 90.42308
$ws.Range("L61").Value = 271.26924
$ws.Range("N61").Value = -701.26924
$ws.Range("H70").Value = 2004
$ws.Range("I70").Value = 2004
$ws.Range("K70").Value = 6012
$ws.Range("M70").Value = -5697
$ws.Range("H73").Value = 2004
$ws.Range("I73").Value = 2004
$ws.Range("K73").Value = 6012
$ws.Range("M73").Value = -4920
$ws.Range("H75").Value = 1503.125
$ws.Range("I75").Value = 1085.5
$ws.Range("J75").Value = 1642.3334
$ws.Range("K75").Value = 3256.5
$ws.Range("L75").Value = 4927.0002
$ws.Range("M75").Value = -2258.5
$ws.Range("N75").Value = -6923.0002
$ws.Range("H78").Value = 1503.125
$ws.Range("I78").Value = 1085.5
$ws.Range("J78").Value = 1642.3334
$ws.Range("K78").Value = 9769.5
$ws.Range("L78").Value = 14781.0006
$ws.Range("M78").Value = -4777.5
$ws.Range("N78").Value = -24765.0006
$ws.Range("H92").Value = 1155.625
$ws.Range("I92").Value = 1265.8182
$ws.Range("K92").Value = 3797.4546
$ws.Range("M92").Value = -2549.4546
$ws.Range("H140").Value = 1542.4166
$ws.Range("I140").Value = 1334.8334
$ws.Range("J140").Value = 1750
$ws.Range("K140").Value = 4004.5002
$ws.Range("L140").Value = 5250
$ws.Range("M140").Value = 1175.4998
$ws.Range("N140").Value = -15610

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 809.4706
$ws.Range("I97").Value = 906.75
$ws.Range("K97").Value = 906.75
$ws.Range("M97").Value = -410.75
$ws.Range("H102").Value = 1332.1111
$ws.Range("I102").Value = 1123.625
$ws.Range("K102").Value = 1123.625
$ws.Range("M102").Value = 498.375
$ws.Range("H107").Value = 1037.4
$ws.Range("I107").Value = 986
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 986
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 934
$ws.Range("N107").Value = -5340
$ws.Range("H135").Value = 86666.336
$ws.Range("J135").Value = 86666.336
$ws.Range("L135").Value = 86666.336
$ws.Range("N135").Value = -96806.336

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2796.4546
$ws.Range("I107").Value = 4824.6
$ws.Range("J107").Value = 1106.3334
$ws.Range("K107").Value = 14473.8
$ws.Range("L107").Value = 3319.0002
$ws.Range("M107").Value = -12553.8
$ws.Range("N107").Value = -7159.0002
$ws.Range("H126").Value = 10004355
$ws.Range("J126").Value = 3375
$ws.Range("L126").Value = 10125
$ws.Range("N126").Value = -15065
$ws.Range("H128").Value = 49715
$ws.Range("J128").Value = 49715
$ws.Range("L128").Value = 49715
$ws.Range("N128").Value = -59675
